$p = $ppt.ActivePresentation

# Slide 4 ("Nacrt (2)"): stretch the dependency-plan picture ("Content
# Placeholder 3", nacrtOdvisnosti.gif) to fill as much of the slide as
# possible (commit: "Raztegnjen nacrt, kolikor se da").
# Target EMU box: off (0, 1500174) ext (9776388 x 4071966).
$slide4 = $p.Slides.Item(4)
$pic = $slide4.Shapes.Item("Content Placeholder 3")
$pic.Left = 0
$pic.Top = 118.12393700787402
$pic.Width = 769.7943427086614
$pic.Height = 320.62724409448816
